$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 130
$ws.Range("B130").Value = 7483247
$ws.Range("F130").Value = "Mushuc Runa"
$ws.Range("G130").Value = "Universidad Catolica del Ecuador"
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = "A"
$ws.Range("K130").Value = 3.25
$ws.Range("L130").Value = 3.2
$ws.Range("M130").Value = 2.25
$ws.Range("N130").Value = 3.5
$ws.Range("O130").Value = 3.25
$ws.Range("P130").Value = 2.1
$ws.Range("Q130").Value = 0.5
$ws.Range("R130").Value = 1.775
$ws.Range("S130").Value = 2.025
$ws.Range("U130").Value = 1.9
$ws.Range("V130").Value = 1.9
$ws.Range("X130").Value = -1
$ws.Range("Y130").Value = 1.1
$ws.Range("AA130").Value = 1.025
$ws.Range("AB130").Value = -1
$ws.Range("AC130").Value = 0.8999999999999999

# Row 131
$ws.Range("B131").Value = 7483081
$ws.Range("F131").Value = "Deportivo Cuenca"
$ws.Range("G131").Value = "El Nacional"
$ws.Range("H131").Value = 1
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = "H"
$ws.Range("K131").Value = 2.75
$ws.Range("L131").Value = 3.25
$ws.Range("M131").Value = 2.55
$ws.Range("N131").Value = 3
$ws.Range("O131").Value = 3.3
$ws.Range("P131").Value = 2.3
$ws.Range("Q131").Value = 0.25
$ws.Range("R131").Value = 1.825
$ws.Range("S131").Value = 1.975
$ws.Range("T131").Value = 2.75
$ws.Range("U131").Value = 2
$ws.Range("V131").Value = 1.8
$ws.Range("W131").Value = 2
$ws.Range("Y131").Value = -1
$ws.Range("Z131").Value = 0.825
$ws.Range("AA131").Value = -1
$ws.Range("AC131").Value = 0.8

# Row 133
$ws.Range("B133").Value = 7483189
$ws.Range("F133").Value = "Independiente del Valle"
$ws.Range("G133").Value = "Orense"
$ws.Range("H133").Value = 2
$ws.Range("I133").Value = 2
$ws.Range("J133").Value = "D"
$ws.Range("K133").Value = 1.4
$ws.Range("L133").Value = 4.75
$ws.Range("M133").Value = 7
$ws.Range("N133").Value = 1.4
$ws.Range("O133").Value = 4.5
$ws.Range("P133").Value = 8
$ws.Range("Q133").Value = -1.25
$ws.Range("R133").Value = 1.875
$ws.Range("S133").Value = 1.925
$ws.Range("T133").Value = 2.5
$ws.Range("U133").Value = 1.925
$ws.Range("V133").Value = 1.875
$ws.Range("W133").Value = -1
$ws.Range("X133").Value = 3.5
$ws.Range("Z133").Value = -1
$ws.Range("AA133").Value = 0.925
$ws.Range("AB133").Value = 0.925
$ws.Range("AC133").Value = -1

# Row 135
$ws.Range("B135").Value = 7483306
$ws.Range("F135").Value = "Tecnico Universitario"
$ws.Range("G135").Value = "Club Atletico Libertad"
$ws.Range("H135").Value = 1
$ws.Range("J135").Value = "D"
$ws.Range("K135").Value = 1.5
$ws.Range("L135").Value = 4.333
$ws.Range("M135").Value = 5.75
$ws.Range("N135").Value = 1.533
$ws.Range("O135").Value = 4.2
$ws.Range("P135").Value = 5.5
$ws.Range("Q135").Value = -1
$ws.Range("R135").Value = 1.925
$ws.Range("S135").Value = 1.875
$ws.Range("T135").Value = 2.25
$ws.Range("U135").Value = 1.8
$ws.Range("V135").Value = 2
$ws.Range("W135").Value = -1
$ws.Range("X135").Value = 3.2
$ws.Range("Z135").Value = -1
$ws.Range("AA135").Value = 0.875
$ws.Range("AB135").Value = -0.5
$ws.Range("AC135").Value = 0.5

# Row 136
$ws.Range("B136").Value = 7483188
$ws.Range("F136").Value = "Gualaceo SC"
$ws.Range("G136").Value = "Emelec"
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 2
$ws.Range("J136").Value = "A"
$ws.Range("K136").Value = 3.6
$ws.Range("L136").Value = 3.3
$ws.Range("M136").Value = 2.05
$ws.Range("N136").Value = 2.6
$ws.Range("O136").Value = 3.25
$ws.Range("P136").Value = 2.75
$ws.Range("Q136").Value = 0
$ws.Range("R136").Value = 1.8
$ws.Range("S136").Value = 2
$ws.Range("T136").Value = 2.5
$ws.Range("U136").Value = 1.975
$ws.Range("V136").Value = 1.825
$ws.Range("X136").Value = -1
$ws.Range("Y136").Value = 1.75
$ws.Range("AA136").Value = 1
$ws.Range("AB136").Value = -1
$ws.Range("AC136").Value = 0.825

# Row 137
$ws.Range("B137").Value = 7482832
$ws.Range("F137").Value = "Barcelona Guayaquil"
$ws.Range("G137").Value = "Guayaquil City"
$ws.Range("H137").Value = 2
$ws.Range("I137").Value = 1
$ws.Range("J137").Value = "H"
$ws.Range("K137").Value = 1.363
$ws.Range("L137").Value = 5
$ws.Range("M137").Value = 7.5
$ws.Range("N137").Value = 1.444
$ws.Range("O137").Value = 4
$ws.Range("P137").Value = 8
$ws.Range("Q137").Value = -1.25
$ws.Range("R137").Value = 2.05
$ws.Range("S137").Value = 1.75
$ws.Range("U137").Value = 1.95
$ws.Range("V137").Value = 1.85
$ws.Range("W137").Value = 0.444
$ws.Range("Y137").Value = -1
$ws.Range("Z137").Value = -0.5
$ws.Range("AA137").Value = 0.375
$ws.Range("AB137").Value = 0.95
$ws.Range("AC137").Value = -1

# Row 139
$ws.Range("B139").Value = 7528859
$ws.Range("F139").Value = "Club Atletico Libertad"
$ws.Range("G139").Value = "Cumbaya FC"
$ws.Range("H139").Value = 3
$ws.Range("I139").Value = 1
$ws.Range("J139").Value = "H"
$ws.Range("K139").Value = 1.727
$ws.Range("M139").Value = 4.333
$ws.Range("N139").Value = 1.4
$ws.Range("O139").Value = 4.2
$ws.Range("P139").Value = 7
$ws.Range("Q139").Value = -1.25
$ws.Range("R139").Value = 2
$ws.Range("S139").Value = 1.8
$ws.Range("U139").Value = 1.95
$ws.Range("V139").Value = 1.85
$ws.Range("W139").Value = 0.3999999999999999
$ws.Range("Y139").Value = -1
$ws.Range("Z139").Value = 1
$ws.Range("AA139").Value = -1
$ws.Range("AB139").Value = 0.95
$ws.Range("AC139").Value = -1

# Row 140
$ws.Range("B140").Value = 7528849
$ws.Range("F140").Value = "Guayaquil City"
$ws.Range("G140").Value = "Gualaceo SC"
$ws.Range("H140").Value = 0
$ws.Range("I140").Value = 2
$ws.Range("J140").Value = "A"
$ws.Range("K140").Value = 1.833
$ws.Range("M140").Value = 3.75
$ws.Range("N140").Value = 2.15
$ws.Range("O140").Value = 3.4
$ws.Range("P140").Value = 3
$ws.Range("Q140").Value = -0.25
$ws.Range("R140").Value = 1.825
$ws.Range("S140").Value = 1.975
$ws.Range("U140").Value = 1.85
$ws.Range("V140").Value = 1.95
$ws.Range("W140").Value = -1
$ws.Range("Y140").Value = 2
$ws.Range("Z140").Value = -1
$ws.Range("AA140").Value = 0.9750000000000001
$ws.Range("AB140").Value = -1
$ws.Range("AC140").Value = 0.95

# Row 142
$ws.Range("B142").Value = 7528858
$ws.Range("F142").Value = "Orense"
$ws.Range("G142").Value = "SD Aucas"
$ws.Range("H142").Value = 1
$ws.Range("I142").Value = 2
$ws.Range("J142").Value = "A"
$ws.Range("K142").Value = 2.2
$ws.Range("L142").Value = 3.2
$ws.Range("M142").Value = 3.2
$ws.Range("N142").Value = 1.95
$ws.Range("O142").Value = 3.2
$ws.Range("P142").Value = 3.8
$ws.Range("Q142").Value = -0.5
$ws.Range("R142").Value = 1.95
$ws.Range("S142").Value = 1.85
$ws.Range("U142").Value = 1.85
$ws.Range("V142").Value = 1.95
$ws.Range("W142").Value = -1
$ws.Range("Y142").Value = 2.8
$ws.Range("Z142").Value = -1
$ws.Range("AA142").Value = 0.8500000000000001
$ws.Range("AB142").Value = 0.8500000000000001

# Row 144
$ws.Range("B144").Value = 7528857
$ws.Range("F144").Value = "Universidad Catolica del Ecuador"
$ws.Range("G144").Value = "Barcelona Guayaquil"
$ws.Range("H144").Value = 0
$ws.Range("I144").Value = 1
$ws.Range("K144").Value = 1.533
$ws.Range("L144").Value = 4
$ws.Range("M144").Value = 5.5
$ws.Range("N144").Value = 1.5
$ws.Range("O144").Value = 4.333
$ws.Range("P144").Value = 5.25
$ws.Range("Q144").Value = -1
$ws.Range("R144").Value = 1.8
$ws.Range("S144").Value = 2
$ws.Range("T144").Value = 3
$ws.Range("U144").Value = 1.975
$ws.Range("V144").Value = 1.825
$ws.Range("Y144").Value = 4.25
$ws.Range("AA144").Value = 1
$ws.Range("AB144").Value = -1
$ws.Range("AC144").Value = 0.825

# Row 145
$ws.Range("B145").Value = 7528848
$ws.Range("F145").Value = "Emelec"
$ws.Range("G145").Value = "Deportivo Cuenca"
$ws.Range("H145").Value = 2
$ws.Range("J145").Value = "H"
$ws.Range("K145").Value = 1.75
$ws.Range("L145").Value = 3.5
$ws.Range("M145").Value = 4.2
$ws.Range("N145").Value = 2.4
$ws.Range("O145").Value = 3.1
$ws.Range("P145").Value = 2.75
$ws.Range("Q145").Value = -0.25
$ws.Range("R145").Value = 2.05
$ws.Range("S145").Value = 1.75
$ws.Range("T145").Value = 2.25
$ws.Range("U145").Value = 1.8
$ws.Range("V145").Value = 2
$ws.Range("W145").Value = 1.4
$ws.Range("Y145").Value = -1
$ws.Range("Z145").Value = 1.05
$ws.Range("AA145").Value = -1
$ws.Range("AB145").Value = 0.8
$ws.Range("AC145").Value = -1

# Row 148
$ws.Range("E148").Value = 45353.625
